# Fruta / hortaliza, semanal
# Insert 3 new rows (one new week of data) above the existing Kiwi block
# at rows 133-135, pushing the previously existing rows 133-141 down to 136-144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 133; existing rows 133:141 shift to 136:144.
$ws.Rows("133:135").Insert()

# Common (constant) values shared by every row in this Kiwi/Vega Monumental block.
$mercadoId   = 11
$mercado     = "Vega Monumental Concepción"
$region      = "Bíobío"
$codreg      = 8
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100101007
$categoria   = "Kiwi"
$variedad    = "Hayward"
$unidad      = "$/bandeja 18 kilos"
$origen      = "Región de O'Higgins"
$kgUnidad    = 18

function Set-KiwiRow($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value2  = $mercadoId
    $ws.Cells.Item($Row, 2).Value2  = $mercado
    $ws.Cells.Item($Row, 3).Value2  = $region
    $ws.Cells.Item($Row, 4).Value2  = $Fecha
    $ws.Cells.Item($Row, 5).Value2  = $codreg
    $ws.Cells.Item($Row, 6).Value2  = $tipo
    $ws.Cells.Item($Row, 7).Value2  = $productoId
    $ws.Cells.Item($Row, 8).Value2  = $producto
    $ws.Cells.Item($Row, 9).Value2  = $categoriaId
    $ws.Cells.Item($Row, 10).Value2 = $categoria
    $ws.Cells.Item($Row, 11).Value2 = $variedad
    $ws.Cells.Item($Row, 12).Value2 = $Calidad
    $ws.Cells.Item($Row, 13).Value2 = $Volumen
    $ws.Cells.Item($Row, 14).Value2 = $PrecioMin
    $ws.Cells.Item($Row, 15).Value2 = $PrecioMax
    $ws.Cells.Item($Row, 16).Value2 = $PrecioProm
    $ws.Cells.Item($Row, 17).Value2 = $unidad
    $ws.Cells.Item($Row, 18).Value2 = $origen
    $ws.Cells.Item($Row, 19).Value2 = $PrecioKg
    $ws.Cells.Item($Row, 20).Value2 = $kgUnidad
}

# New week of data (2022-05-25, serial 44706) inserted at the top of the block.
Set-KiwiRow 133 44706 "Especial" 50 11000 11000 11000 611
Set-KiwiRow 134 44706 "Primera"  50 9000  9000  9000  500
Set-KiwiRow 135 44706 "Segunda"  50 7000  7000  7000  389

Write-Host ("Dimension after edit: {0}" -f $ws.UsedRange.Address())
